$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2:F18").NumberFormat = "@"

$ws.Range("A2").Value = "Identity AI / ML Engineer"
$ws.Range("B2").Value = "Openkyber"
$ws.Range("C2").Value = "AK, US USA"
$ws.Range("D2").Value = 28.9
$ws.Range("E2").Value = "AI Engineer, Generative AI, LangChain, RAG, Hugging Face, FAISS, Pinecone, ChromaDB, Prompt Engineering, TensorFlow"
$ws.Range("F2").Value = "2026-02-24"
$ws.Range("G2").Value = "https://www.indeed.com/viewjob?jk=4fff7807f26ccfa6"

$ws.Range("A3").Value = "Senior Data Scientist - USA Remote"
$ws.Range("B3").Value = "Danaher Diagnostics"
$ws.Range("C3").Value = "New York, NY, US USA"
$ws.Range("D3").Value = 15.6
$ws.Range("E3").Value = "Data Scientist, RAG, LLaMA, Mistral, Hugging Face, Prompt Engineering, MLflow, Git, Snowflake, Databricks"
$ws.Range("F3").Value = "2026-02-24"
$ws.Range("G3").Value = "https://www.indeed.com/viewjob?jk=34dde1de84a8d35d"

$ws.Range("A4").Value = "Senior Data Analyst, Analytics Engineering"
$ws.Range("B4").Value = "Pearl Health"
$ws.Range("C4").Value = "New York, NY, US USA"
$ws.Range("D4").Value = 15.6
$ws.Range("E4").Value = "RAG, Athena, Redshift, BigQuery, CI/CD, Git, Snowflake, BigQuery, Redshift, Python"
$ws.Range("F4").Value = "2026-02-24"
$ws.Range("G4").Value = "https://www.indeed.com/viewjob?jk=ab2e931bbf1e84bf"

$ws.Range("A5").Value = "Full-Stack Senior Software Engineer"
$ws.Range("B5").Value = "nan"
$ws.Range("C5").Value = "Boulder, CO, US USA"
$ws.Range("D5").Value = 14.4
$ws.Range("E5").Value = "RAG, Docker, Kubernetes, CI/CD, GitHub Actions, Terraform, Git, Kafka, PostgreSQL, SQL"
$ws.Range("F5").Value = "2026-02-24"
$ws.Range("G5").Value = "https://www.indeed.com/viewjob?jk=6d1ef7cc52780120"

$ws.Range("A6").Value = "ML Infrastructure Architect"
$ws.Range("B6").Value = "Openkyber"
$ws.Range("C6").Value = "AK, US USA"
$ws.Range("D6").Value = 13.3
$ws.Range("E6").Value = "AI Engineer, Generative AI, LangChain, RAG, Prompt Engineering, TensorFlow, PyTorch, Data Lake, Databricks, Python"
$ws.Range("F6").Value = "2026-02-24"
$ws.Range("G6").Value = "https://www.indeed.com/viewjob?jk=5be5e23fbb3c6051"

$ws.Range("A7").Value = "AI Deployment Engineer"
$ws.Range("B7").Value = "Openkyber"
$ws.Range("C7").Value = "AK, US USA"
$ws.Range("D7").Value = 13.3
$ws.Range("E7").Value = "AI Engineer, RAG, LLaMA, Gemini, Prompt Engineering, S3, Redshift, Terraform, Redshift, Python"
$ws.Range("F7").Value = "2026-02-24"
$ws.Range("G7").Value = "https://www.indeed.com/viewjob?jk=403a6d5b03cee33e"

$ws.Range("A8").Value = "DevOps Engineer"
$ws.Range("B8").Value = "SWAP"
$ws.Range("C8").Value = "US USA"
$ws.Range("D8").Value = 13.3
$ws.Range("E8").Value = "Docker, Kubernetes, CI/CD, GitHub Actions, Terraform, Git, NoSQL, Python, SQL, R"
$ws.Range("F8").Value = "2026-02-24"
$ws.Range("G8").Value = "https://www.indeed.com/viewjob?jk=9db067fbe776dd20"

$ws.Range("A9").Value = "Software Development Engineer in Test"
$ws.Range("B9").Value = "Alteryx"
$ws.Range("C9").Value = "Remote, US USA"
$ws.Range("D9").Value = 12.2
$ws.Range("E9").Value = "RAG, Docker, Kubernetes, CI/CD, Jenkins, Git, Python, SQL, R, Java"
$ws.Range("F9").Value = "2026-02-24"
$ws.Range("G9").Value = "https://www.indeed.com/viewjob?jk=01a62888fe35817c"

$ws.Range("A10").Value = "Senior Software Engineer"
$ws.Range("B10").Value = "Just Appraised"
$ws.Range("C10").Value = "Remote, US USA"
$ws.Range("D10").Value = 12.2
$ws.Range("E10").Value = "RAG, Docker, CI/CD, GitHub Actions, Terraform, Git, PostgreSQL, SQL, R, Java"
$ws.Range("F10").Value = "2026-02-24"
$ws.Range("G10").Value = "https://www.indeed.com/viewjob?jk=b2c46b380ee0c308"

$ws.Range("A11").Value = "AI Data Scientist"
$ws.Range("B11").Value = "The Hartford"
$ws.Range("C11").Value = "Columbus, OH, US USA"
$ws.Range("D11").Value = 11.1
$ws.Range("E11").Value = "Data Scientist, LangChain, RAG, TensorFlow, PyTorch, Git, Matplotlib, Python, SQL, R"
$ws.Range("F11").Value = "2026-02-24"
$ws.Range("G11").Value = "https://www.indeed.com/viewjob?jk=b9aa74a73bfb8b29"

$ws.Range("A12").Value = "AI Data Scientist"
$ws.Range("B12").Value = "The Hartford"
$ws.Range("C12").Value = "Charlotte, NC, US USA"
$ws.Range("D12").Value = 11.1
$ws.Range("E12").Value = "Data Scientist, LangChain, RAG, TensorFlow, PyTorch, Git, Matplotlib, Python, SQL, R"
$ws.Range("F12").Value = "2026-02-24"
$ws.Range("G12").Value = "https://www.indeed.com/viewjob?jk=840834a0721f9428"

$ws.Range("A13").Value = "AI Data Scientist"
$ws.Range("B13").Value = "The Hartford"
$ws.Range("C13").Value = "Hartford, CT, US USA"
$ws.Range("D13").Value = 11.1
$ws.Range("E13").Value = "Data Scientist, LangChain, RAG, TensorFlow, PyTorch, Git, Matplotlib, Python, SQL, R"
$ws.Range("F13").Value = "2026-02-24"
$ws.Range("G13").Value = "https://www.indeed.com/viewjob?jk=c990c19acf0c31b5"

$ws.Range("A14").Value = "AI Data Scientist"
$ws.Range("B14").Value = "The Hartford"
$ws.Range("C14").Value = "Chicago, IL, US USA"
$ws.Range("D14").Value = 11.1
$ws.Range("E14").Value = "Data Scientist, LangChain, RAG, TensorFlow, PyTorch, Git, Matplotlib, Python, SQL, R"
$ws.Range("F14").Value = "2026-02-24"
$ws.Range("G14").Value = "https://www.indeed.com/viewjob?jk=63286a781e5a666e"

$ws.Range("A15").Value = "Sr Data Engineer"
$ws.Range("B15").Value = "AdventHealth Corporate"
$ws.Range("C15").Value = "Altamonte Springs, FL, US USA"
$ws.Range("D15").Value = 11.1
$ws.Range("E15").Value = "RAG, Snowflake, Kafka, Hadoop, Tableau, Python, SQL, R, Java, Optimization"
$ws.Range("F15").Value = "2026-02-24"
$ws.Range("G15").Value = "https://www.indeed.com/viewjob?jk=49fe852730ed983b"

$ws.Range("A16").Value = "AI Full Stack Engineer"
$ws.Range("B16").Value = "Divya Stores"
$ws.Range("C16").Value = "Chicago, IL, US USA"
$ws.Range("D16").Value = 11.1
$ws.Range("E16").Value = "Generative AI, RAG, Gemini, Copilot, FastAPI, AKS, CI/CD, Git, Python, R"
$ws.Range("F16").Value = "2026-02-24"
$ws.Range("G16").Value = "https://www.indeed.com/viewjob?jk=56b238b0f6fd6cb7"

$ws.Range("A17").Value = "Senior Software Engineer - ML Platform"
$ws.Range("B17").Value = "Latitude AI"
$ws.Range("C17").Value = "Pittsburgh, PA, US USA"
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = "RAG, PyTorch, Kubernetes, Terraform, Python, SQL, R, Java, Scala"
$ws.Range("F17").Value = "2026-02-24"
$ws.Range("G17").Value = "https://www.indeed.com/viewjob?jk=4388f870668aad41"

$ws.Range("A18").Value = "Perception Engineer - Data"
$ws.Range("B18").Value = "Forterra"
$ws.Range("C18").Value = "Arlington, VA, US USA"
$ws.Range("D18").Value = 10
$ws.Range("E18").Value = "RAG, TensorFlow, PyTorch, Docker, Kubernetes, CI/CD, Python, R, Optimization"
$ws.Range("F18").Value = "2026-02-24"
$ws.Range("G18").Value = "https://www.indeed.com/viewjob?jk=e4d12d00e8a3f926"
